$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(2, 8).Value = [double]"1.591366224351631e-07"
$ws.Cells.Item(2, 9).Value = [double]"1.591366224351631e-07"
$ws.Cells.Item(2, 12).Value = [double]"45.53646573247122"
$ws.Cells.Item(2, 13).Value = "[30.675961806222297, 60.396969658720145]"
$ws.Cells.Item(2, 14).Value = [double]"1.735213426279358e-07"
$ws.Cells.Item(2, 15).Value = [double]"1.735213426279358e-07"
$ws.Cells.Item(2, 16).Value = [double]"1.515763422452734"
$ws.Cells.Item(2, 17).Value = "[1.1132370364071944, 1.9182898084982734]"
$ws.Cells.Item(2, 18).Value = [double]"1.397336246711234e-09"
$ws.Cells.Item(2, 19).Value = [double]"1.397336246711234e-09"
$ws.Cells.Item(2, 20).Value = [double]"56.66223100531182"
$ws.Cells.Item(2, 21).Value = "[47.026234100624606, 66.29822790999904]"
$ws.Cells.Item(2, 22).Value = [double]"1.998401444325282e-15"
$ws.Cells.Item(2, 23).Value = [double]"1.998401444325282e-15"
$ws.Cells.Item(2, 24).Value = [double]"18.19503503503527"
$ws.Cells.Item(2, 25).Value = [double]"16.65877877877899"
$ws.Cells.Item(2, 26).Value = [double]"19.73129129129154"
$ws.Cells.Item(3, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(3, 8).Value = [double]"8.905145942872394e-08"
$ws.Cells.Item(3, 9).Value = [double]"8.905145942872394e-08"
$ws.Cells.Item(3, 12).Value = [double]"56.74225464502296"
$ws.Cells.Item(3, 13).Value = "[37.40439081516, 76.08011847488592]"
$ws.Cells.Item(3, 14).Value = [double]"4.247852452010648e-07"
$ws.Cells.Item(3, 15).Value = [double]"4.247852452010648e-07"
$ws.Cells.Item(3, 16).Value = [double]"1.691868716347656"
$ws.Cells.Item(3, 17).Value = "[1.2893423303021168, 2.0943951023931957]"
$ws.Cells.Item(3, 18).Value = [double]"7.314993055729246e-11"
$ws.Cells.Item(3, 19).Value = [double]"7.314993055729246e-11"
$ws.Cells.Item(3, 20).Value = [double]"57.52361363964329"
$ws.Cells.Item(3, 21).Value = "[45.82041499010087, 69.2268122891857]"
$ws.Cells.Item(3, 22).Value = [double]"7.125411372044255e-13"
$ws.Cells.Item(3, 23).Value = [double]"7.125411372044255e-13"
$ws.Cells.Item(3, 24).Value = [double]"17.52292292292315"
$ws.Cells.Item(3, 25).Value = [double]"15.98666666666687"
$ws.Cells.Item(3, 26).Value = [double]"19.05917917917943"
$ws.Cells.Item(4, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(4, 8).Value = [double]"2.948714716843881e-08"
$ws.Cells.Item(4, 9).Value = [double]"2.948714716843881e-08"
$ws.Cells.Item(4, 12).Value = [double]"53.36623733456896"
$ws.Cells.Item(4, 13).Value = "[34.78204658543726, 71.95042808370067]"
$ws.Cells.Item(4, 14).Value = [double]"6.533582299628193e-07"
$ws.Cells.Item(4, 15).Value = [double]"6.533582299628193e-07"
$ws.Cells.Item(4, 16).Value = [double]"1.855395060678656"
$ws.Cells.Item(4, 17).Value = "[1.478026573760963, 2.23276354759635]"
$ws.Cells.Item(4, 18).Value = [double]"7.058797990566745e-13"
$ws.Cells.Item(4, 19).Value = [double]"7.058797990566745e-13"
$ws.Cells.Item(4, 20).Value = [double]"54.38513817205492"
$ws.Cells.Item(4, 21).Value = "[43.88763319412297, 64.88264314998688]"
$ws.Cells.Item(4, 22).Value = [double]"1.35003119794419e-13"
$ws.Cells.Item(4, 23).Value = [double]"1.35003119794419e-13"
$ws.Cells.Item(4, 24).Value = [double]"16.89881881881904"
$ws.Cells.Item(4, 25).Value = [double]"15.45857857857878"
$ws.Cells.Item(4, 26).Value = [double]"18.3390590590593"
$ws.Cells.Item(5, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(5, 8).Value = [double]"5.803756364386459e-10"
$ws.Cells.Item(5, 9).Value = [double]"5.803756364386459e-10"
$ws.Cells.Item(5, 12).Value = [double]"59.24526823242874"
$ws.Cells.Item(5, 13).Value = "[43.346193975558705, 75.14434248929877]"
$ws.Cells.Item(5, 14).Value = [double]"1.826634399293425e-09"
$ws.Cells.Item(5, 15).Value = [double]"1.826634399293425e-09"
$ws.Cells.Item(5, 16).Value = [double]"1.654131867655887"
$ws.Cells.Item(5, 17).Value = "[1.3522370781217328, 1.9560266571900415]"
$ws.Cells.Item(5, 18).Value = [double]"2.176037128265307e-14"
$ws.Cells.Item(5, 19).Value = [double]"2.176037128265307e-14"
$ws.Cells.Item(5, 20).Value = [double]"57.66746739849875"
$ws.Cells.Item(5, 21).Value = "[47.64878101450371, 67.68615378249379]"
$ws.Cells.Item(5, 22).Value = [double]"4.218847493575595e-15"
$ws.Cells.Item(5, 23).Value = [double]"4.218847493575595e-15"
$ws.Cells.Item(5, 24).Value = [double]"17.66694694694717"
$ws.Cells.Item(5, 25).Value = [double]"16.51475475475497"
$ws.Cells.Item(5, 26).Value = [double]"18.81913913913938"
$ws.Cells.Item(6, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(6, 8).Value = [double]"1.063813627188992e-07"
$ws.Cells.Item(6, 9).Value = [double]"1.063813627188992e-07"
$ws.Cells.Item(6, 12).Value = [double]"47.31009754195681"
$ws.Cells.Item(6, 13).Value = "[32.08123271068154, 62.538962373232074]"
$ws.Cells.Item(6, 14).Value = [double]"1.295796201450372e-07"
$ws.Cells.Item(6, 15).Value = [double]"1.295796201450372e-07"
$ws.Cells.Item(6, 16).Value = [double]"1.578658170272349"
$ws.Cells.Item(6, 17).Value = "[1.2012896833546556, 1.9560266571900424]"
$ws.Cells.Item(6, 18).Value = [double]"8.346701108052912e-11"
$ws.Cells.Item(6, 19).Value = [double]"8.346701108052912e-11"
$ws.Cells.Item(6, 20).Value = [double]"55.43543997271819"
$ws.Cells.Item(6, 21).Value = "[45.599720636228, 65.27115930920837]"
$ws.Cells.Item(6, 22).Value = [double]"8.43769498715119e-15"
$ws.Cells.Item(6, 23).Value = [double]"8.43769498715119e-15"
$ws.Cells.Item(6, 24).Value = [double]"17.95499499499522"
$ws.Cells.Item(6, 25).Value = [double]"16.51475475475496"
$ws.Cells.Item(6, 26).Value = [double]"19.39523523523549"
$ws.Cells.Item(7, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(7, 8).Value = [double]"1.942553393696755e-06"
$ws.Cells.Item(7, 9).Value = [double]"1.942553393696755e-06"
$ws.Cells.Item(7, 12).Value = [double]"45.75084809287657"
$ws.Cells.Item(7, 13).Value = "[27.095553446698943, 64.40614273905419]"
$ws.Cells.Item(7, 14).Value = [double]"1.125083435526975e-05"
$ws.Cells.Item(7, 15).Value = [double]"1.125083435526975e-05"
$ws.Cells.Item(7, 16).Value = [double]"1.717026615475503"
$ws.Cells.Item(7, 17).Value = "[1.2516054816103486, 2.1824477493406578]"
$ws.Cells.Item(7, 18).Value = [double]"2.354027639484002e-09"
$ws.Cells.Item(7, 19).Value = [double]"2.354027639484002e-09"
$ws.Cells.Item(7, 20).Value = [double]"54.04494574875036"
$ws.Cells.Item(7, 21).Value = "[43.1553814003898, 64.93451009711092]"
$ws.Cells.Item(7, 22).Value = [double]"5.264677582772492e-13"
$ws.Cells.Item(7, 23).Value = [double]"5.264677582772492e-13"
$ws.Cells.Item(7, 24).Value = [double]"17.42690690690713"
$ws.Cells.Item(7, 25).Value = [double]"15.65061061061081"
$ws.Cells.Item(7, 26).Value = [double]"19.20320320320345"
$ws.Cells.Item(8, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(8, 8).Value = [double]"4.897982019969049e-11"
$ws.Cells.Item(8, 9).Value = [double]"4.897982019969049e-11"
$ws.Cells.Item(8, 12).Value = [double]"62.66894467843619"
$ws.Cells.Item(8, 13).Value = "[46.404920243451144, 78.93296911342124]"
$ws.Cells.Item(8, 14).Value = [double]"7.701477233723608e-10"
$ws.Cells.Item(8, 15).Value = [double]"7.701477233723608e-10"
$ws.Cells.Item(8, 16).Value = [double]"1.767342413731195"
$ws.Cells.Item(8, 17).Value = "[1.490605523324887, 2.0440793041375036]"
$ws.Cells.Item(8, 18).Value = [double]"0"
$ws.Cells.Item(8, 19).Value = [double]"0"
$ws.Cells.Item(8, 20).Value = [double]"64.00670641994172"
$ws.Cells.Item(8, 21).Value = "[54.29450252865939, 73.71891031122405]"
$ws.Cells.Item(8, 22).Value = [double]"0"
$ws.Cells.Item(8, 23).Value = [double]"0"
$ws.Cells.Item(8, 24).Value = [double]"17.2348748748751"
$ws.Cells.Item(8, 25).Value = [double]"16.17869869869891"
$ws.Cells.Item(8, 26).Value = [double]"18.29105105105129"
$ws.Cells.Item(9, 6).Value = [double]"23.98000000000031"
$ws.Cells.Item(9, 8).Value = [double]"1.891747053706094e-05"
$ws.Cells.Item(9, 9).Value = [double]"1.891747053706094e-05"
$ws.Cells.Item(9, 12).Value = [double]"40.35455045090148"
$ws.Cells.Item(9, 13).Value = "[21.450320916271828, 59.25877998553113]"
$ws.Cells.Item(9, 14).Value = [double]"9.072625621309172e-05"
$ws.Cells.Item(9, 15).Value = [double]"9.072625621309172e-05"
$ws.Cells.Item(9, 16).Value = [double]"1.754763464167272"
$ws.Cells.Item(9, 17).Value = "[1.2012896833546556, 2.308237244979889]"
$ws.Cells.Item(9, 18).Value = [double]"8.341627233399151e-08"
$ws.Cells.Item(9, 19).Value = [double]"8.341627233399151e-08"
$ws.Cells.Item(9, 20).Value = [double]"59.62888415469533"
$ws.Cells.Item(9, 21).Value = "[48.76662838596681, 70.49113992342384]"
$ws.Cells.Item(9, 22).Value = [double]"2.042810365310288e-14"
$ws.Cells.Item(9, 23).Value = [double]"2.042810365310288e-14"
$ws.Cells.Item(9, 24).Value = [double]"17.2828828828831"
$ws.Cells.Item(9, 25).Value = [double]"15.17053053053072"
$ws.Cells.Item(9, 26).Value = [double]"19.39523523523549"
$ws.Cells.Item(10, 6).Value = [double]"23.66000000000026"
$ws.Cells.Item(10, 8).Value = [double]"2.220354925142054e-06"
$ws.Cells.Item(10, 9).Value = [double]"2.220354925142054e-06"
$ws.Cells.Item(10, 12).Value = [double]"50.32208072777819"
$ws.Cells.Item(10, 13).Value = "[31.91741199811547, 68.7267494574409]"
$ws.Cells.Item(10, 14).Value = [double]"1.673450875872007e-06"
$ws.Cells.Item(10, 15).Value = [double]"1.673450875872007e-06"
$ws.Cells.Item(10, 16).Value = [double]"1.540921321580579"
$ws.Cells.Item(10, 17).Value = "[1.0880791372793475, 1.9937635058818106]"
$ws.Cells.Item(10, 18).Value = [double]"1.680733840991877e-08"
$ws.Cells.Item(10, 19).Value = [double]"1.680733840991877e-08"
$ws.Cells.Item(10, 20).Value = [double]"59.25120336806009"
$ws.Cells.Item(10, 21).Value = "[47.27161426367506, 71.23079247244512]"
$ws.Cells.Item(10, 22).Value = [double]"5.861977570020827e-13"
$ws.Cells.Item(10, 23).Value = [double]"5.861977570020827e-13"
$ws.Cells.Item(10, 24).Value = [double]"17.85749749749769"
$ws.Cells.Item(10, 25).Value = [double]"16.15227227227245"
$ws.Cells.Item(10, 26).Value = [double]"19.56272272272294"
$ws.Cells.Item(11, 6).Value = [double]"23.66000000000026"
$ws.Cells.Item(11, 8).Value = [double]"5.21552678822701e-09"
$ws.Cells.Item(11, 9).Value = [double]"5.21552678822701e-09"
$ws.Cells.Item(11, 12).Value = [double]"50.83556632726996"
$ws.Cells.Item(11, 13).Value = "[33.737566768383516, 67.9335658861564]"
$ws.Cells.Item(11, 14).Value = [double]"3.249680884387374e-07"
$ws.Cells.Item(11, 15).Value = [double]"3.249680884387374e-07"
$ws.Cells.Item(11, 16).Value = [double]"1.956026657190042"
$ws.Cells.Item(11, 17).Value = "[1.591237119836272, 2.320816194543811]"
$ws.Cells.Item(11, 18).Value = [double]"4.418687638008123e-14"
$ws.Cells.Item(11, 19).Value = [double]"4.418687638008123e-14"
$ws.Cells.Item(11, 20).Value = [double]"53.44182735921638"
$ws.Cells.Item(11, 21).Value = "[44.15691184738418, 62.726742871048586]"
$ws.Cells.Item(11, 22).Value = [double]"4.218847493575595e-15"
$ws.Cells.Item(11, 23).Value = [double]"4.218847493575595e-15"
$ws.Cells.Item(11, 24).Value = [double]"16.29437437437455"
$ws.Cells.Item(11, 25).Value = [double]"14.92072072072088"
$ws.Cells.Item(11, 26).Value = [double]"17.66802802802822"
$ws.Cells.Item(12, 6).Value = [double]"23.66000000000026"
$ws.Cells.Item(12, 8).Value = [double]"8.083733682440197e-11"
$ws.Cells.Item(12, 9).Value = [double]"8.083733682440197e-11"
$ws.Cells.Item(12, 12).Value = [double]"64.22854954892125"
$ws.Cells.Item(12, 13).Value = "[48.607444757273214, 79.84965434056929]"
$ws.Cells.Item(12, 14).Value = [double]"1.347741918067413e-10"
$ws.Cells.Item(12, 15).Value = [double]"1.347741918067413e-10"
$ws.Cells.Item(12, 16).Value = [double]"1.603816069400195"
$ws.Cells.Item(12, 17).Value = "[1.3270791789938867, 1.8805529598065034]"
$ws.Cells.Item(12, 18).Value = [double]"3.33066907387547e-15"
$ws.Cells.Item(12, 19).Value = [double]"3.33066907387547e-15"
$ws.Cells.Item(12, 20).Value = [double]"57.32814543481366"
$ws.Cells.Item(12, 21).Value = "[47.27124145141994, 67.38504941820739]"
$ws.Cells.Item(12, 22).Value = [double]"5.773159728050814e-15"
$ws.Cells.Item(12, 23).Value = [double]"5.773159728050814e-15"
$ws.Cells.Item(12, 24).Value = [double]"17.62066066066085"
$ws.Cells.Item(12, 25).Value = [double]"16.57857857857876"
$ws.Cells.Item(12, 26).Value = [double]"18.66274274274295"
$ws.Cells.Item(13, 6).Value = [double]"23.66000000000026"
$ws.Cells.Item(13, 8).Value = [double]"8.274990552781247e-07"
$ws.Cells.Item(13, 9).Value = [double]"8.274990552781247e-07"
$ws.Cells.Item(13, 12).Value = [double]"50.86046598211599"
$ws.Cells.Item(13, 13).Value = "[29.976970890704678, 71.7439610735273]"
$ws.Cells.Item(13, 14).Value = [double]"1.260392966551116e-05"
$ws.Cells.Item(13, 15).Value = [double]"1.260392966551116e-05"
$ws.Cells.Item(13, 16).Value = [double]"1.880552959806503"
$ws.Cells.Item(13, 17).Value = "[1.4402897250691957, 2.320816194543811]"
$ws.Cells.Item(13, 18).Value = [double]"4.645062112729192e-11"
$ws.Cells.Item(13, 19).Value = [double]"4.645062112729192e-11"
$ws.Cells.Item(13, 20).Value = [double]"53.74701852007279"
$ws.Cells.Item(13, 21).Value = "[42.20645164427024, 65.28758539587534]"
$ws.Cells.Item(13, 22).Value = [double]"3.711919660531748e-12"
$ws.Cells.Item(13, 23).Value = [double]"3.711919660531748e-12"
$ws.Cells.Item(13, 24).Value = [double]"16.57857857857876"
$ws.Cells.Item(13, 25).Value = [double]"14.92072072072088"
$ws.Cells.Item(13, 26).Value = [double]"18.23643643643663"
$ws.Cells.Item(14, 6).Value = [double]"23.66000000000026"
$ws.Cells.Item(14, 8).Value = [double]"3.351827257858098e-06"
$ws.Cells.Item(14, 9).Value = [double]"3.351827257858098e-06"
$ws.Cells.Item(14, 12).Value = [double]"52.02754593609544"
$ws.Cells.Item(14, 13).Value = "[28.06538091071744, 75.98971096147343]"
$ws.Cells.Item(14, 14).Value = [double]"7.167369416083424e-05"
$ws.Cells.Item(14, 15).Value = [double]"7.167369416083424e-05"
$ws.Cells.Item(14, 16).Value = [double]"2.132131951084965"
$ws.Cells.Item(14, 17).Value = "[1.6792897667837332, 2.5849741353861964]"
$ws.Cells.Item(14, 18).Value = [double]"2.669642285013651e-12"
$ws.Cells.Item(14, 19).Value = [double]"2.669642285013651e-12"
$ws.Cells.Item(14, 20).Value = [double]"65.08366406982026"
$ws.Cells.Item(14, 21).Value = "[52.388748378011925, 77.77857976162859]"
$ws.Cells.Item(14, 22).Value = [double]"1.887379141862766e-13"
$ws.Cells.Item(14, 23).Value = [double]"1.887379141862766e-13"
$ws.Cells.Item(14, 24).Value = [double]"15.6312312312314"
$ws.Cells.Item(14, 25).Value = [double]"13.92600600600616"
$ws.Cells.Item(14, 26).Value = [double]"17.33645645645665"
